$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("74").Insert()

$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C74").Value = "Los Lagos"
$ws.Range("D74").Value = 44477
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = 100112023
$ws.Range("G74").Value = "Brócoli"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 1050
$ws.Range("N74").Value = "$/unidad"
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 1050
$ws.Range("Q74").Value = 1
$ws.Range("R74").Value = "Hortaliza"
